# Auto-generated Excel COM-interop script applying the CryCompanywiseStockReport_1.xlsx diff.
# Each block sets the literal values for the changed cells in a given row, exactly as
# specified by the target OOXML diff (quantities, values, rates and rolled-up subtotals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38
$ws.Cells.Item(38, 6).Value = 27
$ws.Cells.Item(38, 7).Value = 843.21

# Row 39
$ws.Cells.Item(39, 6).Value = 59
$ws.Cells.Item(39, 7).Value = 4132.36

# Row 44
$ws.Cells.Item(44, 6).Value = 33
$ws.Cells.Item(44, 7).Value = 1164.57

# Row 63
$ws.Cells.Item(63, 2).Value = 30037.41

# Row 151
$ws.Cells.Item(151, 6).Value = 26
$ws.Cells.Item(151, 7).Value = 1286.48

# Row 158
$ws.Cells.Item(158, 6).Value = 106
$ws.Cells.Item(158, 7).Value = 4719.12

# Row 161
$ws.Cells.Item(161, 2).Value = 16996.09

# Row 166
$ws.Cells.Item(166, 2).Value = 53925

# Row 167
$ws.Cells.Item(167, 2).Value = 57756

# Row 179
$ws.Cells.Item(179, 6).Value = 1
$ws.Cells.Item(179, 7).Value = 1108.25

# Row 180
$ws.Cells.Item(180, 6).Value = 5
$ws.Cells.Item(180, 7).Value = 6045

# Row 181
$ws.Cells.Item(181, 2).Value = 59749.1

# Row 232
$ws.Cells.Item(232, 6).Value = 83
$ws.Cells.Item(232, 7).Value = 5719.53

# Row 234
$ws.Cells.Item(234, 2).Value = 6293.28

# Row 283
$ws.Cells.Item(283, 6).Value = 23
$ws.Cells.Item(283, 7).Value = 865.95

# Row 296
$ws.Cells.Item(296, 2).Value = 111480.92

# Row 299
$ws.Cells.Item(299, 6).Value = 75
$ws.Cells.Item(299, 7).Value = 16216.5

# Row 328
$ws.Cells.Item(328, 6).Value = 160
$ws.Cells.Item(328, 7).Value = 8340.799999999999

# Row 329
$ws.Cells.Item(329, 6).Value = 123
$ws.Cells.Item(329, 7).Value = 13771.08

# Row 340
$ws.Cells.Item(340, 6).Value = 55
$ws.Cells.Item(340, 7).Value = 5560.5

# Row 347
$ws.Cells.Item(347, 6).Value = 3
$ws.Cells.Item(347, 7).Value = 177.39

# Row 358
$ws.Cells.Item(358, 6).Value = 118
$ws.Cells.Item(358, 7).Value = 8296.58

# Row 366
$ws.Cells.Item(366, 2).Value = 282332.78

# Row 413
$ws.Cells.Item(413, 2).Value = 58047
$ws.Cells.Item(413, 4).Value = 105.54
$ws.Cells.Item(413, 5).Value = 126.1
$ws.Cells.Item(413, 6).Value = 62
$ws.Cells.Item(413, 7).Value = 6543.48

# Row 414
$ws.Cells.Item(414, 2).Value = 47097
$ws.Cells.Item(414, 4).Value = 112.28
$ws.Cells.Item(414, 5).Value = 134.16
$ws.Cells.Item(414, 6).Value = 15
$ws.Cells.Item(414, 7).Value = 1684.2

# Row 420
$ws.Cells.Item(420, 6).Value = 40
$ws.Cells.Item(420, 7).Value = 1488.4

# Row 424
$ws.Cells.Item(424, 2).Value = 40985.76

# Row 459
$ws.Cells.Item(459, 6).Value = 500
$ws.Cells.Item(459, 7).Value = 6725

# Row 461
$ws.Cells.Item(461, 6).Value = 565
$ws.Cells.Item(461, 7).Value = 7237.65

# Row 462
$ws.Cells.Item(462, 6).Value = 269
$ws.Cells.Item(462, 7).Value = 7074.7

# Row 466
$ws.Cells.Item(466, 6).Value = 351
$ws.Cells.Item(466, 7).Value = 6925.23

# Row 467
$ws.Cells.Item(467, 6).Value = 361
$ws.Cells.Item(467, 7).Value = 2375.38

# Row 470
$ws.Cells.Item(470, 6).Value = 962
$ws.Cells.Item(470, 7).Value = 6329.96

# Row 473
$ws.Cells.Item(473, 6).Value = 417
$ws.Cells.Item(473, 7).Value = 5483.55

# Row 474
$ws.Cells.Item(474, 6).Value = 330
$ws.Cells.Item(474, 7).Value = 8679

# Row 475
$ws.Cells.Item(475, 6).Value = 270
$ws.Cells.Item(475, 7).Value = 4436.1

# Row 477
$ws.Cells.Item(477, 2).Value = 98108.10000000001

# Row 509
$ws.Cells.Item(509, 6).Value = 11
$ws.Cells.Item(509, 7).Value = 5518.7

# Row 510
$ws.Cells.Item(510, 6).Value = 0
$ws.Cells.Item(510, 7).Value = 0

# Row 515
$ws.Cells.Item(515, 2).Value = 18710.57

# Row 525
$ws.Cells.Item(525, 6).Value = 349
$ws.Cells.Item(525, 7).Value = 4498.61

# Row 527
$ws.Cells.Item(527, 6).Value = 368
$ws.Cells.Item(527, 7).Value = 2465.6

# Row 531
$ws.Cells.Item(531, 2).Value = 38526.45

# Row 593
$ws.Cells.Item(593, 6).Value = 106
$ws.Cells.Item(593, 7).Value = 2823.84

# Row 598
$ws.Cells.Item(598, 2).Value = 30519.92

# Row 605
$ws.Cells.Item(605, 6).Value = 25
$ws.Cells.Item(605, 7).Value = 680

# Row 607
$ws.Cells.Item(607, 2).Value = 56247.08

# Row 631
$ws.Cells.Item(631, 6).Value = 69
$ws.Cells.Item(631, 7).Value = 2979.42

# Row 635
$ws.Cells.Item(635, 2).Value = 23690.78

# Row 714
$ws.Cells.Item(714, 6).Value = 314
$ws.Cells.Item(714, 7).Value = 25609.84

# Row 718
$ws.Cells.Item(718, 6).Value = 331
$ws.Cells.Item(718, 7).Value = 43195.5

# Row 724
$ws.Cells.Item(724, 6).Value = 376
$ws.Cells.Item(724, 7).Value = 14017.28

# Row 725
$ws.Cells.Item(725, 6).Value = 87
$ws.Cells.Item(725, 7).Value = 7889.16

# Row 728
$ws.Cells.Item(728, 6).Value = 1
$ws.Cells.Item(728, 7).Value = 53.73

# Row 731
$ws.Cells.Item(731, 6).Value = 646
$ws.Cells.Item(731, 7).Value = 87216.46000000001

# Row 732
$ws.Cells.Item(732, 6).Value = 51
$ws.Cells.Item(732, 7).Value = 1908.42

# Row 733
$ws.Cells.Item(733, 6).Value = 712
$ws.Cells.Item(733, 7).Value = 85945.52

# Row 735
$ws.Cells.Item(735, 2).Value = 315296.9

# Row 758
$ws.Cells.Item(758, 6).Value = 16
$ws.Cells.Item(758, 7).Value = 598.4

# Row 761
$ws.Cells.Item(761, 2).Value = 2319.42

# Row 807
$ws.Cells.Item(807, 6).Value = 3401
$ws.Cells.Item(807, 7).Value = 554737.11

# Row 815
$ws.Cells.Item(815, 2).Value = 699566.83

# Row 821
$ws.Cells.Item(821, 2).Value = 3195051.88

# Row 822
$ws.Cells.Item(822, 2).Value = 3195051.88
